$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") '28.340.97'
Set-TextValue $ws.Range("E2") '  +0.72%  '
Set-TextValue $ws.Range("D3") '1.804.74'
Set-TextValue $ws.Range("E3") '  -0.76%  '
Set-TextValue $ws.Range("D4") '1.003'
Set-TextValue $ws.Range("E4") '  +0.23%  '
Set-TextValue $ws.Range("D5") '326.88'
Set-TextValue $ws.Range("E5") '  -3.29%  '
Set-TextValue $ws.Range("D6") '0.9984'
Set-TextValue $ws.Range("E6") '  +0.00%  '
Set-TextValue $ws.Range("D7") '0.4451'
Set-TextValue $ws.Range("E7") '  +4.41%  '
Set-TextValue $ws.Range("D8") '0.3756'
Set-TextValue $ws.Range("E8") '  +5.90%  '
Set-TextValue $ws.Range("D9") '44.75'
Set-TextValue $ws.Range("E9") '  -1.95%  '
Set-TextValue $ws.Range("D10") '1.150'
Set-TextValue $ws.Range("E10") '  -1.45%  '
Set-TextValue $ws.Range("D11") '0.07519'
Set-TextValue $ws.Range("E11") '  -0.17%  '
Set-TextValue $ws.Range("D12") '22.55'
Set-TextValue $ws.Range("E12") '  -2.10%  '
Set-TextValue $ws.Range("D13") '1.002'
Set-TextValue $ws.Range("E13") '  +0.28%  '
Set-TextValue $ws.Range("D14") '7.700'
Set-TextValue $ws.Range("E14") '  +5.06%  '
Set-TextValue $ws.Range("D15") '6.300'
Set-TextValue $ws.Range("E15") '  -0.36%  '
Set-TextValue $ws.Range("D16") '1.805.76'
Set-TextValue $ws.Range("E16") '  -0.42%  '
Set-TextValue $ws.Range("D17") '0.00001093'
Set-TextValue $ws.Range("E17") '  -0.37%  '
Set-TextValue $ws.Range("D18") '0.06786'
Set-TextValue $ws.Range("E18") '  +1.38%  '
Set-TextValue $ws.Range("D19") '80.80'
Set-TextValue $ws.Range("E19") '  -2.31%  '
Set-TextValue $ws.Range("D20") '0.9993'
Set-TextValue $ws.Range("E20") '  +0.12%  '
Set-TextValue $ws.Range("E21") '  +0.04%  '
Set-TextValue $ws.Range("D22") '6.328'
Set-TextValue $ws.Range("E22") '  -1.19%  '
Set-TextValue $ws.Range("D23") '28.365.27'
Set-TextValue $ws.Range("E23") '  +0.70%  '
Set-TextValue $ws.Range("E24") '  -1.02%  '
Set-TextValue $ws.Range("D25") '2.407'
Set-TextValue $ws.Range("E25") '  +0.11%  '
Set-TextValue $ws.Range("D26") '20.49'
Set-TextValue $ws.Range("E26") '  -2.03%  '
Set-TextValue $ws.Range("D27") '153.16'
Set-TextValue $ws.Range("E27") '  -1.78%  '
Set-TextValue $ws.Range("D28") '2.354'
Set-TextValue $ws.Range("E28") '  -6.54%  '
Set-TextValue $ws.Range("D29") '2.007.08'
Set-TextValue $ws.Range("E29") '  -0.71%  '
Set-TextValue $ws.Range("D30") '132.74'
Set-TextValue $ws.Range("E30") '  -0.55%  '
Set-TextValue $ws.Range("D31") '1.250'
Set-TextValue $ws.Range("E31") '  -6.02%  '
Set-TextValue $ws.Range("D32") '4.014'
Set-TextValue $ws.Range("E32") '  -1.46%  '
Set-TextValue $ws.Range("D33") '5.843'
Set-TextValue $ws.Range("E33") '  -3.45%  '
Set-TextValue $ws.Range("D34") '0.09323'
Set-TextValue $ws.Range("E34") '  +1.10%  '
Set-TextValue $ws.Range("D35") '0.2280'
Set-TextValue $ws.Range("E35") '  +5.16%  '
Set-TextValue $ws.Range("E36") '  -2.70%  '
Set-TextValue $ws.Range("D37") '0.06350'
Set-TextValue $ws.Range("E37") '  -0.19%  '
Set-TextValue $ws.Range("D38") '0.02320'
Set-TextValue $ws.Range("E38") '  -1.72%  '
Set-TextValue $ws.Range("D39") '5.170'
Set-TextValue $ws.Range("E39") '  -2.09%  '
Set-TextValue $ws.Range("D40") '0.6580'
Set-TextValue $ws.Range("E40") '  -2.07%  '
Set-TextValue $ws.Range("D41") '1.209'
Set-TextValue $ws.Range("E41") '  -1.22%  '
Set-TextValue $ws.Range("D42") '1.460'
Set-TextValue $ws.Range("E42") '  -3.60%  '
Set-TextValue $ws.Range("D43") '8.152'
Set-TextValue $ws.Range("E43") '  -0.73%  '
Set-TextValue $ws.Range("D44") '0.9985'
Set-TextValue $ws.Range("E44") '  +0.02%  '
Set-TextValue $ws.Range("D45") '13.91'
Set-TextValue $ws.Range("E45") '  -3.00%  '
Set-TextValue $ws.Range("E46") '  -2.17%  '
Set-TextValue $ws.Range("D47") '3.793'
Set-TextValue $ws.Range("E47") '  -2.24%  '
Set-TextValue $ws.Range("D48") '128.45'
Set-TextValue $ws.Range("E48") '  -0.54%  '
Set-TextValue $ws.Range("D49") '2.034'
Set-TextValue $ws.Range("E49") '  -1.64%  '
Set-TextValue $ws.Range("D50") '0.07104'
Set-TextValue $ws.Range("E50") '  -0.27%  '
Set-TextValue $ws.Range("D51") '1.159'
Set-TextValue $ws.Range("E51") '  -2.53%  '
